# Updated symbol list (Price / Volume(1h) columns) to reflect the latest
# coinranking.com snapshot. Both columns are stored as plain text in the
# sheet (e.g. "330.85", "0.04%"), so force a text number format before
# writing each value - otherwise Excel auto-converts the numeric-looking
# / percent-looking strings into actual numbers and we lose the exact
# textual formatting (trailing zeros, "%", sign, etc.).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "330.85"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.04%"

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "41.56"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "0.95%"

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.688"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.24%"

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08394"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "4.03%"

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "8.816"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.71%"

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "2.002"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-1.52%"

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "4.480"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-0.98%"

# Row 9
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-0.84%"

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9250"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "0.44%"

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1281"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "1.80%"

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.1975"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "1.51%"

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09446"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "1.17%"

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.03952"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "6.81%"

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.1062"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.95%"

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001299"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-0.10%"

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.006110"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-2.21%"

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.423"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "1.80%"

# Row 19
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.73%"

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.976"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "8.14%"

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1363"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-3.82%"

# Row 22
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-5.49%"

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04402"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.89%"

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001245"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-1.31%"

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004378"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "2.04%"

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001192"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-4.05%"

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0003993"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "0.01%"

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02831"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-1.01%"

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05527"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "1.04%"

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007952"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "2.16%"

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1438"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "1.45%"

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.008968"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-9.92%"

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002094"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-6.12%"

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.01178"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-0.33%"

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006940"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "2.50%"

# Row 47
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.10%"

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003303"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "9.43%"

# Row 49
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-0.07%"

# Row 50
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.10%"

# Row 51
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.10%"
